# Atualizacao Visual e documentacao
#
# Appends a new block of forecast rows (hora_relatorio = "23:00:00") to
# Sheet1, mirroring the existing 20:00:00 / 21:00:00 blocks already present
# in rows 2-25 / 26-49. This extends the used range from A1:G49 to A1:G73
# and introduces one new shared-string value ("23:00:00").
#
# Columns B ("day") and F ("data_relatorio") hold ISO-like date strings
# (e.g. "2024-07-04") that must stay literal text -- format those columns
# as Text up front so Excel doesn't auto-convert them to date serials.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1,   "2024-07-02", 1,   1,   "o3",   "2024-07-04", "23:00:00"),
    @(3,   "2024-07-03", 24,  1,   "o3",   "2024-07-04", "23:00:00"),
    @(2,   "2024-07-04", 22,  1,   "o3",   "2024-07-04", "23:00:00"),
    @(2,   "2024-07-05", 24,  1,   "o3",   "2024-07-04", "23:00:00"),
    @(2,   "2024-07-06", 38,  1,   "o3",   "2024-07-04", "23:00:00"),
    @(3,   "2024-07-07", 22,  1,   "o3",   "2024-07-04", "23:00:00"),
    @(3,   "2024-07-08", 7,   1,   "o3",   "2024-07-04", "23:00:00"),
    @(3,   "2024-07-09", 8,   1,   "o3",   "2024-07-04", "23:00:00"),
    @(178, "2024-07-02", 178, 178, "pm25", "2024-07-04", "23:00:00"),
    @(146, "2024-07-03", 175, 68,  "pm25", "2024-07-04", "23:00:00"),
    @(162, "2024-07-04", 184, 72,  "pm25", "2024-07-04", "23:00:00"),
    @(137, "2024-07-05", 190, 59,  "pm25", "2024-07-04", "23:00:00"),
    @(107, "2024-07-06", 162, 45,  "pm25", "2024-07-04", "23:00:00"),
    @(106, "2024-07-07", 157, 44,  "pm25", "2024-07-04", "23:00:00"),
    @(50,  "2024-07-08", 164, 22,  "pm25", "2024-07-04", "23:00:00"),
    @(41,  "2024-07-09", 53,  16,  "pm25", "2024-07-04", "23:00:00"),
    @(103, "2024-07-02", 103, 103, "pm10", "2024-07-04", "23:00:00"),
    @(70,  "2024-07-03", 100, 27,  "pm10", "2024-07-04", "23:00:00"),
    @(87,  "2024-07-04", 115, 30,  "pm10", "2024-07-04", "23:00:00"),
    @(70,  "2024-07-05", 121, 25,  "pm10", "2024-07-04", "23:00:00"),
    @(58,  "2024-07-06", 84,  30,  "pm10", "2024-07-04", "23:00:00"),
    @(52,  "2024-07-07", 75,  15,  "pm10", "2024-07-04", "23:00:00"),
    @(34,  "2024-07-08", 83,  14,  "pm10", "2024-07-04", "23:00:00"),
    @(20,  "2024-07-09", 25,  10,  "pm10", "2024-07-04", "23:00:00")
)

$startRow = 50
$endRow = $startRow + $data.Count - 1

$ws.Range("B${startRow}:B${endRow}").NumberFormat = "@"
$ws.Range("F${startRow}:F${endRow}").NumberFormat = "@"

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $values = $data[$i]
    for ($c = 0; $c -lt $values.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $values[$c]
    }
}
